$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: lowercase / rename columns ---
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- Data rows: carbon/CED values moved from column E into column D,
#     column E zeroed out, column F populated with climate-change impact ---
$rowData = @(
    @(2, 1.0, 0, 0.000027882653),
    @(3, 0, 0, 0),
    @(4, 29.8, 0, 0.00083090307),
    @(5, 27.0, 0, 0.00075283164),
    @(6, 0, 0, 0),
    @(7, 0, 0, 0),
    @(8, 273.0000000000001, 0, 0.0076119644),
    @(9, 0.006, 0, 0.00000016729592),
    @(10, 0.437, 0, 0.00001218472),
    @(11, 0, 0, 0),
    @(12, 0.02, 0, 0.00000055765307),
    @(13, 25200.0, 0, 0.70264287),
    @(14, 1530.0, 0, 0.04266046),
    @(15, 5810.0, 0, 0.16199822),
    @(16, 771.0000000000001, 0, 0.021497526),
    @(17, 0.501, 0, 0.000013969209),
    @(18, 0, 0, 0),
    @(19, 0, 0, 0),
    @(20, 0, 0, 0),
    @(21, 0, 0, 0),
    @(22, 0, 0, 0),
    @(23, 0, 0, 0),
    @(24, 0, 0, 0),
    @(25, 0, 0, 0),
    @(26, 0, 0, 0),
    @(27, 0, 0, 0),
    @(28, 0, 0, 0),
    @(29, 0, 0, 0),
    @(30, 0, 0, 0),
    @(31, 0, 0, 0),
    @(32, 20.6, 0, 0.00057438266),
    @(33, 0, 0, 0),
    @(34, 0, 0, 0),
    @(35, 0, 0, 0),
    @(36, 0, 0, 0),
    @(37, 0, 0, 0),
    @(38, 0, 0, 0),
    @(39, 0, 0, 0),
    @(40, 0, 0, 0),
    @(41, 0, 0, 0),
    @(42, 0, 0, 0),
    @(43, 0, 0, 0),
    @(44, 0, 0, 0),
    @(45, 0, 0, 0),
    @(46, 6.34, 0, 0.00017677602),
    @(47, 0, 0, 0),
    @(48, 0.044, 0, 0.0000012268368),
    @(49, 0, 0, 0),
    @(50, 0, 0, 0),
    @(51, 0, 0, 0),
    @(52, 0, 0, 0)
)

foreach ($item in $rowData) {
    $r = $item[0]
    $ws.Cells.Item($r, 4).Value = $item[1]
    $ws.Cells.Item($r, 5).Value = $item[2]
    $ws.Cells.Item($r, 6).Value = $item[3]
}

# --- Header cell comments ---
$excel.UserName = "Data Processor"
$ws.Range("A1").AddComment("Data type: Categorical (text)")
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)")
$ws.Range("C1").AddComment("Data type: Categorical (text)")
$ws.Range("D1").AddComment("Data type: Carbon footprint")
$ws.Range("E1").AddComment("Data type: Cumulative energy demand")
$ws.Range("F1").AddComment("Data type: Climate change impact")
$ws.Range("G1").AddComment("Data type: Categorical (text)")
